$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confusion matrix values (rows 2-6, columns C/B/D/E/F) to match
# the refreshed evaluation results referenced in the commit message.
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 8
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 2

$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 16

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 21
$ws.Range("F5").Value = 11

$ws.Range("E6").Value = 9
